$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.02
$ws.Range("G2").Value = 2.04
$ws.Range("H2").Value = 3.6
$ws.Range("I2").Value = 3.65
$ws.Range("K2").Value = 4.5
$ws.Range("N2").Value = 6.8
$ws.Range("O2").Value = 1.15
$ws.Range("P2").Value = 2.94
$ws.Range("Q2").Value = 1.5
$ws.Range("R2").Value = 1.82
$ws.Range("S2").Value = 2.2
$ws.Range("T2").Value = 1.5
$ws.Range("U2").Value = 2.9
$ws.Range("V2").Value = 1.38
$ws.Range("W2").Value = 1.96
$ws.Range("X2").Value = 30
$ws.Range("Y2").Value = 24
$ws.Range("AA2").Value = 65
$ws.Range("AB2").Value = 16.5
$ws.Range("AF2").Value = 17
$ws.Range("AM2").Value = 50
$ws.Range("AN2").Value = 8.199999999999999
$ws.Range("AO2").Value = 19.5
$ws.Range("F3").Value = 1.89
$ws.Range("G3").Value = 2.22
$ws.Range("H3").Value = 3.45
$ws.Range("I3").Value = 4.7
$ws.Range("J3").Value = 3.45
$ws.Range("K3").Value = 4.4
$ws.Range("L3").Value = 1.34
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 3.6
$ws.Range("O3").Value = 1.25
$ws.Range("P3").Value = 2
$ws.Range("Q3").Value = 1.71
$ws.Range("R3").Value = 1.39
$ws.Range("S3").Value = 2.62
$ws.Range("U3").Value = 2.12
$ws.Range("V3").Value = 1.28
$ws.Range("W3").Value = 1.82
$ws.Range("X3").Value = 19.5
$ws.Range("Y3").Value = 18
$ws.Range("Z3").Value = 34
$ws.Range("AA3").Value = 85
$ws.Range("AB3").Value = 12
$ws.Range("AC3").Value = 9.800000000000001
$ws.Range("AD3").Value = 18
$ws.Range("AE3").Value = 50
$ws.Range("AF3").Value = 15.5
$ws.Range("AG3").Value = 12
$ws.Range("AH3").Value = 19
$ws.Range("AI3").Value = 55
$ws.Range("AJ3").Value = 27
$ws.Range("AK3").Value = 23
$ws.Range("AL3").Value = 36
$ws.Range("AM3").Value = 100
$ws.Range("AN3").Value = 14
$ws.Range("AO3").Value = 44
$ws.Range("F4").Value = 5.3
$ws.Range("G4").Value = 6.2
$ws.Range("H4").Value = 1.6
$ws.Range("J4").Value = 4.5
$ws.Range("K4").Value = 5
$ws.Range("N4").Value = 5.7
$ws.Range("F5").Value = 2.18
$ws.Range("G5").Value = 2.42
$ws.Range("J5").Value = 2.84
$ws.Range("K5").Value = 3.3
$ws.Range("T5").Value = 2.26
$ws.Range("W5").Value = 1.7
$ws.Range("AB5").Value = 7.6
$ws.Range("H6").Value = 9.4
$ws.Range("N6").Value = 3.05
$ws.Range("O6").Value = 1.38
$ws.Range("R6").Value = 1.25
$ws.Range("S6").Value = 4.1
$ws.Range("G7").Value = 1.15
$ws.Range("J7").Value = 8.800000000000001
$ws.Range("N7").Value = 2.82
$ws.Range("P7").Value = 2.82
$ws.Range("W7").Value = 6.6
$ws.Range("H8").Value = 10
$ws.Range("I8").Value = 13.5
$ws.Range("J8").Value = 5.6
$ws.Range("R8").Value = 1.69
$ws.Range("U8").Value = 1.99
$ws.Range("V8").Value = 1.08
$ws.Range("G9").Value = 3.85
$ws.Range("N9").Value = 3.8
$ws.Range("P9").Value = 2
$ws.Range("Q9").Value = 1.88
$ws.Range("U9").Value = 2.14
$ws.Range("F10").Value = 2.36
$ws.Range("H10").Value = 2.96
$ws.Range("V10").Value = 1.47
$ws.Range("G11").Value = 2.18
$ws.Range("H11").Value = 3.25
$ws.Range("I11").Value = 3.7
$ws.Range("V11").Value = 1.37
$ws.Range("W11").Value = 1.84
$ws.Range("Z11").Value = 34
$ws.Range("G12").Value = 2.64
$ws.Range("W12").Value = 1.63
$ws.Range("G13").Value = 8
$ws.Range("I13").Value = 1.5
$ws.Range("L13").Value = 1.21
$ws.Range("S13").Value = 2.36
$ws.Range("W13").Value = 1.14
$ws.Range("X13").Value = 32
$ws.Range("AG13").Value = 32
$ws.Range("G14").Value = 2.36
$ws.Range("J14").Value = 3.6
$ws.Range("Q14").Value = 1.8
$ws.Range("W14").Value = 1.74
$ws.Range("F15").Value = 1.22
$ws.Range("G15").Value = 1.25
$ws.Range("Q15").Value = 1.34
$ws.Range("S15").Value = 1.92
$ws.Range("V15").Value = 1.06
$ws.Range("AN15").Value = 3.3
$ws.Range("G16").Value = 2.28
$ws.Range("H16").Value = 3.35
$ws.Range("I16").Value = 4.1
$ws.Range("J16").Value = 3.6
$ws.Range("K16").Value = 4.1
$ws.Range("L16").Value = 1.29
$ws.Range("N16").Value = 4.2
$ws.Range("P16").Value = 2.12
$ws.Range("Q16").Value = 1.73
$ws.Range("R16").Value = 1.43
$ws.Range("S16").Value = 2.92
$ws.Range("U16").Value = 2.28
$ws.Range("Y16").Value = 18.5
$ws.Range("AK16").Value = 980
$ws.Range("AM16").Value = 80
$ws.Range("AN16").Value = 16
$ws.Range("K17").Value = 4.6
$ws.Range("H18").Value = 2.9
$ws.Range("L18").Value = 1.42
$ws.Range("W18").Value = 1.56
$ws.Range("K19").Value = 4.6
$ws.Range("L19").Value = 1.24
$ws.Range("W19").Value = 1.22
$ws.Range("AI19").Value = 1000
$ws.Range("F20").Value = 6
$ws.Range("H20").Value = 1.5
$ws.Range("I20").Value = 1.51
$ws.Range("J20").Value = 5.2
$ws.Range("L20").Value = 1.19
$ws.Range("M20").Value = 1.02
$ws.Range("Q20").Value = 1.41
$ws.Range("S20").Value = 2
$ws.Range("V20").Value = 2.96
$ws.Range("AC20").Value = 16
$ws.Range("H21").Value = 1.72
$ws.Range("X21").Value = 19.5
$ws.Range("Z21").Value = 11.5
$ws.Range("F22").Value = 2.48
$ws.Range("G22").Value = 2.56
$ws.Range("H22").Value = 3.45
$ws.Range("Q22").Value = 2.52
$ws.Range("F24").Value = 1.2
$ws.Range("G24").Value = 1.23
$ws.Range("H24").Value = 12.5
$ws.Range("I24").Value = 15
$ws.Range("J24").Value = 8.199999999999999
$ws.Range("K24").Value = 9.800000000000001
$ws.Range("L24").Value = 1.16
$ws.Range("N24").Value = 10
$ws.Range("P24").Value = 3.9
$ws.Range("R24").Value = 2.14
$ws.Range("S24").Value = 1.71
$ws.Range("T24").Value = 1.69
$ws.Range("U24").Value = 2.18
$ws.Range("V24").Value = 1.07
$ws.Range("W24").Value = 5.1
$ws.Range("X24").Value = 1000
$ws.Range("Y24").Value = 1000
$ws.Range("Z24").Value = 190
$ws.Range("AA24").Value = 540
$ws.Range("AC24").Value = 22
$ws.Range("AD24").Value = 50
$ws.Range("AE24").Value = 180
$ws.Range("AG24").Value = 13
$ws.Range("AH24").Value = 32
$ws.Range("AI24").Value = 140
$ws.Range("AM24").Value = 130
$ws.Range("AN24").Value = 2.92
$ws.Range("F25").Value = 2.6
$ws.Range("H25").Value = 2.82
$ws.Range("I25").Value = 2.84
$ws.Range("P25").Value = 2.4
$ws.Range("R25").Value = 1.56
$ws.Range("S25").Value = 2.7
$ws.Range("U25").Value = 2.68
$ws.Range("Y25").Value = 15
$ws.Range("AD25").Value = 12.5
$ws.Range("AE25").Value = 26
$ws.Range("AF25").Value = 19
$ws.Range("AI25").Value = 30
$ws.Range("AN25").Value = 16.5
